$wb = $excel.ActiveWorkbook

# ---------- Sheet "Rushing" ----------
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2 - M.Ryan
$rushing.Range("C2").Value = 2
$rushing.Range("D2").Value = 6
$rushing.Range("E2").Value = 8
$rushing.Range("F2").Value = 2

# Row 4 - M.Davis
$rushing.Range("C4").Value = 70
$rushing.Range("D4").Value = 24
$rushing.Range("F4").Value = 13

# Row 5 - C.Patterson
$rushing.Range("C5").Value = 72
$rushing.Range("D5").Value = 27
$rushing.Range("E5").Value = 8
$rushing.Range("F5").Value = 17

# Row 7 - Q.Ollison
$rushing.Range("C7").Value = 8

# Row 8 - was K.Smith, now becomes O.Zaccheaus (new player inserted here,
# pushing K.Smith down to the new row 9)
$rushing.Range("B8").Value = "O.Zaccheaus"
$rushing.Range("C8").Value = 0
$rushing.Range("D8").Value = 0
$rushing.Range("F8").Value = 0

# Row 9 - new row, K.Smith's original stats move down here unchanged
$rushing.Range("A8").Copy()
$rushing.Range("A9").PasteSpecial(-4122)
$rushing.Range("A9").Value = 7
$rushing.Range("B9").Value = "K.Smith"
$rushing.Range("C9").Value = 1
$rushing.Range("D9").Value = 1
$rushing.Range("E9").Value = 1
$rushing.Range("F9").Value = 1

# ---------- Sheet "Receiving" ----------
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2 - M.Davis
$receiving.Range("C2").Value = 44
$receiving.Range("D2").Value = 33
$receiving.Range("E2").Value = 1
$receiving.Range("F2").Value = 1

# Row 3 - C.Patterson
$receiving.Range("C3").Value = 47
$receiving.Range("D3").Value = 36
$receiving.Range("E3").Value = 10
$receiving.Range("G3").Value = 10

# Row 5 - Q.Ollison
$receiving.Range("C5").Value = 4
$receiving.Range("D5").Value = 3

# Row 6 - R.Gage
$receiving.Range("C6").Value = 56
$receiving.Range("D6").Value = 49
$receiving.Range("G6").Value = 7
$receiving.Range("H6").Value = 4

# Row 7 - O.Zaccheaus
$receiving.Range("C7").Value = 34
$receiving.Range("D7").Value = 19
$receiving.Range("E7").Value = 5
$receiving.Range("F7").Value = 2
$receiving.Range("G7").Value = 6

# Row 9 - T.Sharpe
$receiving.Range("C9").Value = 27
$receiving.Range("D9").Value = 23
$receiving.Range("E9").Value = 4
$receiving.Range("G9").Value = 2

# Row 10 - K.Pitts
$receiving.Range("C10").Value = 60
$receiving.Range("D10").Value = 37
$receiving.Range("E10").Value = 21
$receiving.Range("F10").Value = 12

# Row 14 - K.Smith
$receiving.Range("C14").Value = 5
$receiving.Range("D14").Value = 5
